# ooutput update 2025 august
#
# Updates the canonical-URL / ValueSet-URL host (migrated from the old
# github.com/RicardoLSantos/shorthand location to the new 2rdoc.pt IG site)
# and refreshes the IG generation timestamp. Also re-applies the "best fit"
# column widths on the Elements sheet that Excel recalculated once the
# shorter hostname changed the overall autofit metrics for this workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet — canonical URL + generation Date
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/advanced-vital-signs-context"
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# ---------------------------------------------------------------------
# 2. Elements sheet — ValueSet binding URL
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/advanced-vital-signs-context-vs"

# ---------------------------------------------------------------------
# 3. Elements sheet — refreshed best-fit column widths.
#    (columns 3, 4, 31, 32 and 33 stay hidden; width is re-applied then
#    the hidden flag is restored since setting ColumnWidth clears it)
# ---------------------------------------------------------------------
$widths = @{
    1  = 15.666666666666666
    2  = 15.666666666666666
    3  = 9.0
    4  = 6.166666666666667
    5  = 4.5
    6  = 3.1666666666666665
    7  = 3.5
    8  = 11.833333333333334
    9  = 9.666666666666666
    11 = 13.5
    15 = 11.5
    20 = 7.0
    21 = 12.833333333333334
    22 = 13.166666666666666
    23 = 14.166666666666666
    24 = 13.833333333333334
    25 = 16.166666666666668
    26 = 66.5
    27 = 4.166666666666667
    28 = 17.166666666666668
    29 = 33.666666666666664
    30 = 12.666666666666666
    31 = 10.5
    32 = 14.166666666666666
    33 = 7.333333333333333
    34 = 7.666666666666667
    37 = 18.666666666666668
}
$hiddenCols = @(3, 4, 31, 32, 33)

foreach ($col in $widths.Keys) {
    $elements.Columns.Item($col).ColumnWidth = $widths[$col]
}
foreach ($col in $hiddenCols) {
    $elements.Columns.Item($col).Hidden = $true
}
